# refactor with sew to implement code chunks
#
# The sheet gains a new first row holding a "code chunk" label
# ("1:3; 4:6") rendered in a monospace (Courier New) font, which pushes
# the pre-existing values (1..6, previously on rows 1-3 and 5-7) down by
# two rows (now rows 3-5 and 7-9), keeping the original blank-row gap
# between the two groups of three numbers intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new label row by inserting two rows above the
# current row 1; this shifts every existing row (and its values) down
# by two, reproducing rows 3,4,5 and 7,8,9 from the original 1,2,3 and
# 5,6,7.
$ws.Rows("1:2").Insert()

# Write the label text into the freshly inserted A1 and give it the
# monospace font used for code chunks.
$ws.Range("A1").Value = "1:3; 4:6"
$ws.Range("A1").Font.Name = "Courier New"
